# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to the refreshed values from the regenerated gh-pages output.

$wb = $excel.ActiveWorkbook

# row number (r) -> new value for column F
$updates = @{
    4  = 870
    5  = 38
    6  = 336
    7  = 10562
    8  = 155
    9  = 88
    10 = 3
    11 = 80
    12 = 137
    13 = 136
    16 = 37
    19 = 297
    20 = 985
    22 = 102
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
